# "Subclustering of fresh samples"
# Re-number the Cluster ID column on Fresh_v_control to reflect the new
# sub-clustering of the fresh samples, refresh the (no-fill) formatting
# that Excel re-applied to the re-sorted rows, and move the active
# sheet/selection from Fresh_mixed back to Fresh_v_control.

$wb = $excel.ActiveWorkbook

$wsFresh   = $wb.Worksheets.Item("Fresh_v_control")
$wsMixed   = $wb.Worksheets.Item("Fresh_mixed")

# --- Updated Cluster IDs (column A) for Fresh_v_control --------------------
# Annotation (column B) / UMAP (column C) per row stay exactly as they were;
# only the numeric cluster id in column A was renumbered for each row.
$newClusterIds = @{
    2  = 8
    3  = 3
    4  = 6
    5  = 1
    6  = 7
    7  = 10
    8  = 13
    9  = 12
    10 = 11
    11 = 16
    12 = 15
    13 = 5
    14 = 14
    15 = 4
    16 = 9
    17 = 2
}

foreach ($row in $newClusterIds.Keys) {
    $wsFresh.Cells.Item($row, 1).Value = $newClusterIds[$row]
}

# Re-apply "no fill" formatting across the data rows (A2:C17) -- this is the
# formatting pass Excel performs when the rows are refreshed after the
# re-sort, and is what produced the extra (no-colour) cell styles in the
# saved workbook.
$wsFresh.Range("A2:C17").Interior.ColorIndex = -4142

# --- Sheet / selection activation ------------------------------------------
# Fresh_mixed was the active tab with B11 selected; the edit moves focus to
# Fresh_v_control with B13 selected, and leaves A2:B10 selected on
# Fresh_mixed.
$wsMixed.Activate()
$wsMixed.Range("A2:B10").Select()

$wsFresh.Activate()
$wsFresh.Range("B13").Select()
